# Update res_bus/vm_pu.xlsx results for the "case with 380 kV" run.
# Slack bus voltage setpoint (column B) changes from 1.05 to 1.02 p.u.,
# and the downstream bus voltage magnitudes (columns C:F, I:N) are refreshed
# with the corresponding load-flow results for every data row (r2:r25).
# Columns A and G (and the untouched H) keep their original values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Worksheet column indices for B,C,D,E,F,I,J,K,L,M,N
$cols = @(2,3,4,5,6,9,10,11,12,13,14)

# New vm_pu values per row, in the same column order as $cols
$data = @{
    2 = @(1.02, 1.049736449159357, 1.055841343850541, 1.046878816966142, 1.064686919665635, 1.046835850159026, 1.054773316231763, 1.058580693466262, 1.049643018994857, 1.06740218526035, 1.056271215325002)
    3 = @(1.02, 1.051042882300646, 1.056603859327626, 1.048002975686036, 1.06566561341653, 1.047170413351613, 1.055727403704518, 1.059156860045573, 1.050578098213877, 1.068195723151123, 1.057226657711404)
    4 = @(1.02, 1.051887844116897, 1.05709717228694, 1.048730286175898, 1.066299027273642, 1.047385606066681, 1.056343889264246, 1.05952892963672, 1.051182477912135, 1.068708708712492, 1.057844018751338)
    5 = @(1.02, 1.052242976180588, 1.057304540527321, 1.049036026420399, 1.066565346861778, 1.047475764191129, 1.056602852608992, 1.059685168695449, 1.051436398177482, 1.068924251884223, 1.058103349853738)
    6 = @(1.02, 1.052302599231679, 1.05733935736207, 1.049087360367023, 1.066610065004973, 1.047490884035625, 1.056646321596152, 1.059711391431484, 1.05147902311543, 1.068960435746229, 1.058146880571845)
    7 = @(1.02, 1.051892589758411, 1.057099943232252, 1.048734371572686, 1.066302585721167, 1.047386811977163, 1.056347350357022, 1.059531018013709, 1.051185871437604, 1.068711589265906, 1.057847484759263)
    8 = @(1.02, 1.050178046427726, 1.056099057211121, 1.047258752469124, 1.065017646407193, 1.046949185138353, 1.055095936354762, 1.058775566834813, 1.049959174976482, 1.067670466153778, 1.056594293605581)
    9 = @(1.02, 1.047153699652173, 1.054334714374455, 1.04465771398433, 1.062754422098638, 1.046168115949116, 1.052884020358555, 1.057438619511804, 1.04779230545248, 1.065832124397619, 1.054379236434949)
    10 = @(1.02, 1.045135195535931, 1.05315804511616, 1.04292303195617, 1.061246249172809, 1.045640708700605, 1.051404747865531, 1.056543439668, 1.046344076350257, 1.064604006293335, 1.05289786320531)
    11 = @(1.02, 1.044260579612514, 1.052648428612865, 1.042171714041723, 1.060593338339958, 1.045410740966798, 1.050763074854185, 1.056154891419481, 1.045716087811886, 1.064071602330194, 1.052255278944687)
    12 = @(1.02, 1.043935615898249, 1.052459117479165, 1.041892610385447, 1.060350838131951, 1.045325080147615, 1.050524555586465, 1.056010427155375, 1.045482688190267, 1.063873749680407, 1.052016420952248)
    13 = @(1.02, 1.044005325960183, 1.052499726126924, 1.041952480495018, 1.060402854359866, 1.045343465580509, 1.050575726660609, 1.056041421580385, 1.04553275942254, 1.063916194018039, 1.052067664695187)
    14 = @(1.02, 1.044233719917428, 1.052632780440931, 1.042148643886732, 1.060573292794798, 1.045403665127927, 1.050743362326717, 1.056142952820177, 1.045696797727732, 1.064055249693291, 1.052235538423169)
    15 = @(1.02, 1.044374428614415, 1.05271475735294, 1.042269502485146, 1.060678308159136, 1.045440724159579, 1.050846625154091, 1.056205490982334, 1.04579784896109, 1.064140913987653, 1.052338947895603)
    16 = @(1.02, 1.045193228158871, 1.053191864330232, 1.042972890269394, 1.061289583554986, 1.045655937201555, 1.051447309409473, 1.056569206710021, 1.046385734826097, 1.064639327040152, 1.052940485191527)
    17 = @(1.02, 1.045706678805308, 1.053491111085306, 1.043414054572201, 1.061673056853731, 1.045790506662656, 1.051823796385655, 1.056797106778408, 1.046754258765557, 1.064951801592998, 1.053317506822399)
    18 = @(1.02, 1.046006109176578, 1.053665645965051, 1.043671360153814, 1.061896743710496, 1.045868844717715, 1.052043285202424, 1.05692994738737, 1.046969126260422, 1.065134002881871, 1.053537307338462)
    19 = @(1.02, 1.046108197689438, 1.053725156041155, 1.043759091708381, 1.061973017438935, 1.04589552986233, 1.052118106685415, 1.056975227411516, 1.047042375952331, 1.065196118689863, 1.053612235076542)
    20 = @(1.02, 1.045651596303694, 1.053459005848329, 1.043366723706845, 1.061631912421591, 1.045776084570682, 1.051783414257385, 1.056772664534519, 1.046714728575734, 1.064918282226269, 1.053277067346874)
    21 = @(1.02, 1.044166466170783, 1.05259359972923, 1.04209087952405, 1.060523102379522, 1.045385944507572, 1.050694002605068, 1.056113058273395, 1.045648496326666, 1.064014303867684, 1.052186108605054)
    22 = @(1.02, 1.043232169354876, 1.052049386820806, 1.041288526699273, 1.05982606450445, 1.045139255817421, 1.05000804220412, 1.055697526551501, 1.044977321725067, 1.063445391858995, 1.051499174061655)
    23 = @(1.02, 1.043727509880899, 1.052337893714449, 1.041713886849001, 1.060195566833178, 1.045270162334917, 1.050371778743966, 1.05591788482978, 1.045333199874165, 1.063747034819923, 1.05186342714911)
    24 = @(1.02, 1.045676485890269, 1.053473512856165, 1.043388110543477, 1.061650503778278, 1.04578260176923, 1.051801661541749, 1.056783709209518, 1.046732590833552, 1.064933428370247, 1.053295340544476)
    25 = @(1.02, 1.047935952665291, 1.054790916851586, 1.045330252202871, 1.063339404144522, 1.046371218915493, 1.053456666984637, 1.057784935230403, 1.04835312877055, 1.066307827835777, 1.054952696284923)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $rowVals[$i]
    }
}
